# Clarify status report for unit tests:
#  - Column D ("PASS/FAIL") formulas now distinguish an errored "Expected"
#    value (column B) from a genuine PASS/FAIL comparison by reporting
#    "ERROR" when B is an error, instead of letting the comparison itself
#    surface the raw #NUM!/#VALUE! error.
#  - A few "Actual" values in column B were refreshed to their latest
#    recalculated results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D3 holds its own (non-shared) formula.
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'

# D4:D36 share formula si="0", anchored at D4; re-enter the whole block so the
# shared-formula group is rewritten consistently.
$ws.Range("D4:D36").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'

# Refreshed "Actual" values (column B) for the Abcd Atm Vol Curve tests.
$ws.Range("B24").Value = 0.16276901888733139
$ws.Range("B25").Value = 0.20553128316863267
$ws.Range("B27").Value = 816.38970366714
$ws.Range("B31").Value = 1.0563094653822682
